$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date (column C) from 45192 to 45202 for all data rows (2-350)
$ws.Range("C2:C350").Value = 45202
